# "Generate Report for Archive"
# The localization status report is regenerated: the handoff status text
# changes from "Ready for handoff" to "In Translation" everywhere it is
# shown (the Overview sheet's per-language status columns, and the
# "Status" column on each per-language detail sheet). Excel's column
# autofit then narrows the affected status columns to match the new
# (shorter) text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- Per-language detail sheets: "Status" column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

# --- Mirror Excel's column-autofit side effect on the status columns ---
# (ColumnWidth is expressed in characters and gets snapped to a pixel
# grid by the host, so this lands on the closest representable width to
# the narrower autofit result rather than an exact value.)
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
